# Scheduled runner update: refresh computed profit figures (columns H-N)
# for a set of recipe rows across the job sheets (ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, WVR). A few rows also gain/lose a trailing cell (M/N) as the
# profit/loss math flips sign and a column becomes (ir)relevant.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 170.16667
$ws.Range("I11").Value = 170.16667
$ws.Range("K11").Value = 170.16667
$ws.Range("M11").Value = -30.16667000000001

$ws.Range("H32").Value = 4497.5
$ws.Range("I32").Value = 4996.6665
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 4996.6665
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -4670.6665
$ws.Range("N32").Value = -3652

$ws.Range("H51").Value = 2958
$ws.Range("I51").Value = 1001
$ws.Range("J51").Value = 3447.25
$ws.Range("K51").Value = 1001
$ws.Range("L51").Value = 3447.25
$ws.Range("M51").Value = -517
$ws.Range("N51").Value = -4415.25

$ws.Range("H76").Value = 5567.1665
$ws.Range("I76").Value = 6250.75
$ws.Range("J76").Value = 4200
$ws.Range("K76").Value = 6250.75
$ws.Range("L76").Value = 4200
$ws.Range("M76").Value = -5935.75
$ws.Range("N76").Value = -4830

$ws.Range("H79").Value = 5567.1665
$ws.Range("I79").Value = 6250.75
$ws.Range("J79").Value = 4200
$ws.Range("K79").Value = 6250.75
$ws.Range("L79").Value = 4200
$ws.Range("M79").Value = -5158.75
$ws.Range("N79").Value = -6384

$ws.Range("H137").Value = 1135.9362
$ws.Range("I137").Value = 880.1
$ws.Range("J137").Value = 2597.8572
$ws.Range("K137").Value = 2640.3
$ws.Range("L137").Value = 7793.571599999999
$ws.Range("M137").Value = -90.30000000000018
$ws.Range("N137").Value = -12893.5716

$ws.Range("H138").Value = 2255.8557
$ws.Range("I138").Value = 1342.7455
$ws.Range("J138").Value = 3451.5952
$ws.Range("K138").Value = 4028.2365
$ws.Range("L138").Value = 10354.7856
$ws.Range("M138").Value = 1111.7635
$ws.Range("N138").Value = -20634.7856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3759.087
$ws.Range("I61").Value = 3041.3125
$ws.Range("K61").Value = 3041.3125
$ws.Range("M61").Value = -2829.3125

$ws.Range("H132").Value = 4328.0713
$ws.Range("I132").Value = 5636.759
$ws.Range("J132").Value = 2922.4443
$ws.Range("K132").Value = 16910.277
$ws.Range("L132").Value = 8767.332900000001
$ws.Range("M132").Value = -14380.277
$ws.Range("N132").Value = -13827.3329

$ws.Range("H136").Value = 3759.087
$ws.Range("I136").Value = 3041.3125
$ws.Range("K136").Value = 9123.9375
$ws.Range("M136").Value = -6573.9375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2834.6943
$ws.Range("I134").Value = 2905.0435
$ws.Range("J134").Value = 2710.2307
$ws.Range("K134").Value = 8715.130500000001
$ws.Range("L134").Value = 8130.6921
$ws.Range("M134").Value = -6180.130500000001
$ws.Range("N134").Value = -13200.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2374.1177
$ws.Range("I31").Value = 1841.9706
$ws.Range("J31").Value = 3438.4119
$ws.Range("K31").Value = 1841.9706
$ws.Range("L31").Value = 3438.4119
$ws.Range("M31").Value = -1546.9706
$ws.Range("N31").Value = -4028.4119

$ws.Range("H34").Value = 2374.1177
$ws.Range("I34").Value = 1841.9706
$ws.Range("J34").Value = 3438.4119
$ws.Range("K34").Value = 1841.9706
$ws.Range("L34").Value = 3438.4119
$ws.Range("M34").Value = -1639.9706
$ws.Range("N34").Value = -3842.4119

$ws.Range("H58").Value = 1373518.9
$ws.Range("I58").Value = 1611791.1
$ws.Range("J58").Value = 3453.25
$ws.Range("K58").Value = 1611791.1
$ws.Range("L58").Value = 3453.25
$ws.Range("M58").Value = -1611588.1
$ws.Range("N58").Value = -3859.25

$ws.Range("H132").Value = 357185.06
$ws.Range("I132").Value = 501699.88
$ws.Range("J132").Value = 2466.818
$ws.Range("K132").Value = 1505099.64
$ws.Range("L132").Value = 7400.454000000001
$ws.Range("M132").Value = -1502569.64
$ws.Range("N132").Value = -12460.454

$ws.Range("H134").Value = 1261.3606
$ws.Range("I134").Value = 1025.4166
$ws.Range("J134").Value = 2132.5386
$ws.Range("K134").Value = 3076.2498
$ws.Range("L134").Value = 6397.6158
$ws.Range("M134").Value = -541.2498000000001
$ws.Range("N134").Value = -11467.6158

$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 1373518.9
$ws.Range("I136").Value = 1611791.1
$ws.Range("J136").Value = 3453.25
$ws.Range("K136").Value = 4835373.300000001
$ws.Range("L136").Value = 10359.75
$ws.Range("M136").Value = -4832823.300000001
$ws.Range("N136").Value = -15459.75

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 86514
$ws.Range("J138").Value = 86514
$ws.Range("L138").Value = 86514
$ws.Range("N138").Value = -96794

$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4208.6
$ws.Range("I137").Value = 3515
$ws.Range("J137").Value = 4671
$ws.Range("K137").Value = 10545
$ws.Range("L137").Value = 14013
$ws.Range("M137").Value = -5445
$ws.Range("N137").Value = -24213

$ws.Range("H138").Value = 2559.6875
$ws.Range("J138").Value = 3551.6
$ws.Range("L138").Value = 10654.8
$ws.Range("N138").Value = -20934.8

$ws.Range("H139").Value = 2471.8823
$ws.Range("I139").Value = 2098.75
$ws.Range("J139").Value = 2803.5557
$ws.Range("K139").Value = 6296.25
$ws.Range("L139").Value = 8410.667099999999
$ws.Range("M139").Value = -1156.25
$ws.Range("N139").Value = -18690.6671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1845.8163
$ws.Range("I132").Value = 1304.375
$ws.Range("J132").Value = 2365.6
$ws.Range("K132").Value = 3913.125
$ws.Range("L132").Value = 7096.799999999999
$ws.Range("M132").Value = -1383.125
$ws.Range("N132").Value = -12156.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2900.6667
$ws.Range("I68").Value = 1567.3334
$ws.Range("J68").Value = 3345.111
$ws.Range("K68").Value = 1567.3334
$ws.Range("L68").Value = 3345.111
$ws.Range("M68").Value = -818.3334
$ws.Range("N68").Value = -4843.111

$ws.Range("H71").Value = 2900.6667
$ws.Range("I71").Value = 1567.3334
$ws.Range("J71").Value = 3345.111
$ws.Range("K71").Value = 7836.666999999999
$ws.Range("L71").Value = 16725.555
$ws.Range("M71").Value = -4092.666999999999
$ws.Range("N71").Value = -24213.555

$ws.Range("H132").Value = 4249.846
$ws.Range("I132").Value = 3772.6365
$ws.Range("J132").Value = 6874.5
$ws.Range("K132").Value = 11317.9095
$ws.Range("L132").Value = 20623.5
$ws.Range("M132").Value = -8787.9095
$ws.Range("N132").Value = -25683.5

$ws.Range("H134").Value = 40000
$ws.Range("J134").Value = 40000
$ws.Range("L134").Value = 40000
$ws.Range("N134").Value = -50140

$ws.Range("H135").Value = 59424.5
$ws.Range("J135").Value = 59424.5
$ws.Range("L135").Value = 59424.5
$ws.Range("N135").Value = -69564.5

$ws.Range("H136").Value = 2000.9153
$ws.Range("I136").Value = 1480.98
$ws.Range("J136").Value = 4889.4443
$ws.Range("K136").Value = 4442.940000000001
$ws.Range("L136").Value = 14668.3329
$ws.Range("M136").Value = -1892.940000000001
$ws.Range("N136").Value = -19768.3329

$ws.Range("H137").Value = 59325
$ws.Range("I137").Value = 48000
$ws.Range("J137").Value = 63100
$ws.Range("K137").Value = 48000
$ws.Range("L137").Value = 63100
$ws.Range("M137").Value = -42900
$ws.Range("N137").Value = -73300

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1171.341
$ws.Range("I132").Value = 760.43243
$ws.Range("J132").Value = 3343.2856
$ws.Range("K132").Value = 2281.29729
$ws.Range("L132").Value = 10029.8568
$ws.Range("M132").Value = 248.70271
$ws.Range("N132").Value = -15089.8568

$ws.Range("H136").Value = 1530.3939
$ws.Range("I136").Value = 1603.5862
$ws.Range("J136").Value = 999.75
$ws.Range("K136").Value = 4810.7586
$ws.Range("L136").Value = 2999.25
$ws.Range("M136").Value = -2260.7586
$ws.Range("N136").Value = -8099.25

$ws.Range("H139").Value = 53610
$ws.Range("I139").Value = 15000
$ws.Range("J139").Value = 61332
$ws.Range("K139").Value = 15000
$ws.Range("L139").Value = 61332
$ws.Range("M139").Value = -9860
$ws.Range("N139").Value = -71612
